# Daily attendance processing - 2025-10-31 19:17:06
# Reverse the order of the comma-separated "Recorded By" entries in column G
# whenever "System" is the first entry of a multi-entry list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text

    if ($text -ne $null -and $text -ne "") {
        $parts = $text.Split(", ")
        $n = $parts.Length

        if ($n -gt 1 -and $parts[0] -eq "System") {
            $reversed = ""
            for ($i = $n - 1; $i -ge 0; $i--) {
                if ($reversed -ne "") {
                    $reversed = $reversed + ", "
                }
                $reversed = $reversed + $parts[$i]
            }
            $cell.Value = $reversed
        }
    }
}
